{"js": "// Update the two-digit multiplication problems in the practice-sheet\n// table. Each entry below maps the original \"A\u00d7B=\" expression (as it\n// appears in before.docx) to its replacement, per the target diff.\nconst replacements = [\n  [\"70\u00d712=\", \"80\u00d726=\"],\n  [\"72\u00d755=\", \"75\u00d743=\"],\n  [\"14\u00d762=\", \"87\u00d755=\"],\n  [\"48\u00d721=\", \"56\u00d714=\"],\n  [\"46\u00d715=\", \"33\u00d760=\"],\n  [\"87\u00d794=\", \"75\u00d714=\"],\n  [\"29\u00d715=\", \"88\u00d723=\"],\n  [\"85\u00d777=\", \"75\u00d728=\"],\n  [\"87\u00d769=\", \"24\u00d762=\"],\n  [\"74\u00d757=\", \"13\u00d714=\"],\n  [\"16\u00d759=\", \"18\u00d758=\"],\n  [\"13\u00d726=\", \"54\u00d754=\"],\n  [\"64\u00d796=\", \"67\u00d738=\"],\n  [\"47\u00d718=\", \"87\u00d713=\"],\n  [\"11\u00d749=\", \"39\u00d769=\"],\n  [\"84\u00d790=\", \"66\u00d748=\"],\n  [\"62\u00d784=\", \"83\u00d752=\"],\n  [\"12\u00d742=\", \"88\u00d755=\"],\n  [\"21\u00d778=\", \"38\u00d757=\"],\n  [\"24\u00d750=\", \"60\u00d757=\"],\n  [\"50\u00d761=\", \"41\u00d750=\"],\n  [\"36\u00d784=\", \"99\u00d743=\"],\n  [\"94\u00d795=\", \"51\u00d742=\"],\n  [\"16\u00d730=\", \"47\u00d771=\"],\n  [\"74\u00d758=\", \"18\u00d795=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems in the practice-sheet\n# table. Each pair below maps the original \"A\u00d7B=\" expression (as it\n# appears in before.docx) to its replacement, per the target diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"70\u00d712=\", \"80\u00d726=\"),\n    @(\"72\u00d755=\", \"75\u00d743=\"),\n    @(\"14\u00d762=\", \"87\u00d755=\"),\n    @(\"48\u00d721=\", \"56\u00d714=\"),\n    @(\"46\u00d715=\", \"33\u00d760=\"),\n    @(\"87\u00d794=\", \"75\u00d714=\"),\n    @(\"29\u00d715=\", \"88\u00d723=\"),\n    @(\"85\u00d777=\", \"75\u00d728=\"),\n    @(\"87\u00d769=\", \"24\u00d762=\"),\n    @(\"74\u00d757=\", \"13\u00d714=\"),\n    @(\"16\u00d759=\", \"18\u00d758=\"),\n    @(\"13\u00d726=\", \"54\u00d754=\"),\n    @(\"64\u00d796=\", \"67\u00d738=\"),\n    @(\"47\u00d718=\", \"87\u00d713=\"),\n    @(\"11\u00d749=\", \"39\u00d769=\"),\n    @(\"84\u00d790=\", \"66\u00d748=\"),\n    @(\"62\u00d784=\", \"83\u00d752=\"),\n    @(\"12\u00d742=\", \"88\u00d755=\"),\n    @(\"21\u00d778=\", \"38\u00d757=\"),\n    @(\"24\u00d750=\", \"60\u00d757=\"),\n    @(\"50\u00d761=\", \"41\u00d750=\"),\n    @(\"36\u00d784=\", \"99\u00d743=\"),\n    @(\"94\u00d795=\", \"51\u00d742=\"),\n    @(\"16\u00d730=\", \"47\u00d771=\"),\n    @(\"74\u00d758=\", \"18\u00d795=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        Write-Output \"WARNING: replacement not found for $oldText\"\n    }\n}\n"}
